$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - cardholder first name
$ws.Range("C2").Value = "Hartmut"

# Row 3 - card number + surname.
# The card number is a 16-digit string; it must stay a text value (as in
# the source file) rather than become a genuine number, otherwise Excel's
# General number display only keeps ~15 significant digits and the last
# two digits would visibly corrupt (...427075 -> ...427080). Force the
# cell to Text format first so the digit string round-trips exactly.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Row 5 - opening balance date label
$ws.Range("D5").Value = "KONTOSTAND AM 21.08.2024"

# Row 6 - transaction 1
$ws.Range("B6").Value = "24.08."
$ws.Range("C6").Value = "25.08."
$ws.Range("D6").Value = "EBAY MKTPLC EU MULGUV"
$ws.Range("E6").Value = "137,65-"

# Row 7 - transaction 2
$ws.Range("B7").Value = "27.08."
$ws.Range("C7").Value = "28.08."
$ws.Range("D7").Value = "KARTENZAHLUNG SHELL TANKSTELLE"
$ws.Range("E7").Value = "63,09-"

# Row 8 - transaction 3
$ws.Range("B8").Value = "28.08."
$ws.Range("C8").Value = "29.08."
$ws.Range("D8").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E8").Value = "25,05-"

# Row 9 - fourth transaction removed entirely; cells go back to blank.
# Restore the blank-row look used elsewhere (rows 10/11): clear the
# text cells and blank out the amount cell, centering it like the
# existing "empty slot" style used for this column.
$ws.Range("B9:D9").Value = ""
$ws.Range("E9").Value = ""
$ws.Range("E9").HorizontalAlignment = -4108   # xlCenter
$ws.Range("E9").VerticalAlignment = -4108     # xlCenter
$ws.Range("E9").WrapText = $true

# Row 12 - closing balance date + amount
$ws.Range("D12").Value = "KONTOSTAND AM 31.08.2024"
$ws.Range("E12").Value = "225,79-"

# Row 13 - next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 07.09.2024"
